$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: iteration headers shift right by two columns worth of labels ---
# F7 was "Iteration #3" -> becomes "Iteration #2" (new Actual sub-column for Iteration #2)
$ws.Range("F7").Value = "Iteration #2"
# G7 was "Transition" -> becomes "Iteration #3" (old Iteration #3 block moves to G:H)
$ws.Range("G7").Value = "Iteration #3"
# New cells: H7 Iteration #3 (Actual), I7/J7 Transition (Estimated/Actual)
$ws.Range("H7").Value = "Iteration #3"
$ws.Range("I7").Value = "Transition"
$ws.Range("J7").Value = "Transition"

# --- Row 8: due-date subheaders ---
# G8 "Due: Apr-24" is removed entirely (Transition now has no specific due date)
$ws.Range("G8").Value = ""
# New H8 due date placeholder for the relocated Transition column
$ws.Range("H8").Value = "Due:"

# --- Row 9: Estimated / Actual sub-headers for the two new column-pairs ---
$ws.Range("G9").Value = "Estimated"
$ws.Range("H9").Value = "Actual"
$ws.Range("I9").Value = "Estimated"
$ws.Range("J9").Value = "Actual"

# --- New "Actual" work-effort values for Iteration #2 (column F, rows 23-30) ---
$ws.Range("F23").Value = 15
$ws.Range("F24").Value = 15
$ws.Range("F25").Value = 5.75
$ws.Range("F26").Value = 5.75
$ws.Range("F27").Value = 2.25
$ws.Range("F28").Value = 7.25
$ws.Range("F29").Value = 7.25
$ws.Range("F30").Value = 4.25

# --- Totals section formulas now reference the summary row instead of re-summing the raw rows ---
$ws.Range("C38").Formula = "=SUM(C36 + E36)"
$ws.Range("C39").Formula = "=SUM(D36 + F36)"

# --- Sheet view bookkeeping to match the saved workbook state ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("J24").Select()

$wb.Application.Calculate()
